# Update the "St. Francis (IL)_B" team-specific transition matrix with the
# newly-computed probabilities (more simulated games -> non-zero transition
# rates for several starting states).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2777777777777778
$ws.Range("C2").Value = 0.3888888888888889
$ws.Range("J2").Value = 0.05555555555555555
$ws.Range("P2").Value = 0.2222222222222222
$ws.Range("S2").Value = 0.05555555555555555

$ws.Range("P3").Value = 0.7142857142857143
$ws.Range("S3").Value = 0.2857142857142857

$ws.Range("B6").Value = 0.125
$ws.Range("F6").Value = 0.25
$ws.Range("J6").Value = 0.125
$ws.Range("Q6").Value = 0.1666666666666667
$ws.Range("R6").Value = 0.04166666666666666
$ws.Range("S6").Value = 0.2916666666666667

$ws.Range("B7").Value = 0.1
$ws.Range("J7").Value = 0.2
$ws.Range("Q7").Value = 0.1
$ws.Range("S7").Value = 0.6

$ws.Range("B8").Value = 0.0303030303030303
$ws.Range("F8").Value = 0.0303030303030303
$ws.Range("J8").Value = 0.1818181818181818
$ws.Range("Q8").Value = 0.1818181818181818
$ws.Range("R8").Value = 0.09090909090909091
$ws.Range("S8").Value = 0.4848484848484849

$ws.Range("B9").Value = 0.05555555555555555
$ws.Range("J9").Value = 0.2222222222222222
$ws.Range("O9").Value = 0.05555555555555555
$ws.Range("Q9").Value = 0.3333333333333333
$ws.Range("R9").Value = 0.05555555555555555
$ws.Range("S9").Value = 0.2777777777777778

$ws.Range("B10").Value = 0.0958904109589041
$ws.Range("F10").Value = 0.1232876712328767
$ws.Range("J10").Value = 0.0684931506849315
$ws.Range("O10").Value = 0.0136986301369863
$ws.Range("Q10").Value = 0.1780821917808219
$ws.Range("R10").Value = 0.0410958904109589
$ws.Range("S10").Value = 0.4794520547945205

$ws.Range("G11").Value = 0.2222222222222222
$ws.Range("J11").Value = 0.05555555555555555
$ws.Range("K11").Value = 0.2222222222222222
$ws.Range("L11").Value = 0.4444444444444444
$ws.Range("S11").Value = 0.05555555555555555

$ws.Range("G12").Value = 0.5
$ws.Range("J12").Value = 0.375
$ws.Range("S12").Value = 0.125

$ws.Range("G13").Value = 0.5
$ws.Range("J13").Value = 0.25
$ws.Range("S13").Value = 0.25

$ws.Range("H15").Value = 0.2857142857142857
$ws.Range("J15").Value = 0.2857142857142857
$ws.Range("K15").Value = 0.07142857142857142
$ws.Range("M15").Value = 0.07142857142857142
$ws.Range("S15").Value = 0.2857142857142857

$ws.Range("I16").Value = 0.125
$ws.Range("J16").Value = 0.5
$ws.Range("O16").Value = 0.125
$ws.Range("S16").Value = 0.25

$ws.Range("F17").Value = 0.06451612903225806
$ws.Range("H17").Value = 0.1935483870967742
$ws.Range("I17").Value = 0.1290322580645161
$ws.Range("J17").Value = 0.3548387096774194
$ws.Range("K17").Value = 0.09677419354838709
$ws.Range("M17").Value = 0.06451612903225806
$ws.Range("S17").Value = 0.09677419354838709

$ws.Range("F18").Value = 0.125
$ws.Range("H18").Value = 0.125
$ws.Range("I18").Value = 0.25
$ws.Range("J18").Value = 0.25
$ws.Range("M18").Value = 0.125
$ws.Range("S18").Value = 0.125

$ws.Range("F19").Value = 0.02150537634408602
$ws.Range("H19").Value = 0.2365591397849462
$ws.Range("I19").Value = 0.1290322580645161
$ws.Range("J19").Value = 0.2903225806451613
$ws.Range("K19").Value = 0.1075268817204301
$ws.Range("M19").Value = 0.01075268817204301
$ws.Range("O19").Value = 0.1075268817204301
$ws.Range("S19").Value = 0.09677419354838709
